# Edit: append two trailing spaces to the first paragraph's existing
# text, then add a new, red-colored run with the "(This is a change -
# Version for branch alternate)" annotation right after it, inside the
# same paragraph.

$d = $word.ActiveDocument

# 1) Add two trailing spaces after "This is a Microsoft word document."
$find = $d.Content.Find
$find.Execute(
    "This is a Microsoft word document.",  # FindText
    $true,                                  # MatchCase
    $false,                                 # MatchWholeWord
    $false,                                 # MatchWildcards
    $false,                                 # MatchSoundsLike
    $false,                                 # MatchAllWordForms
    $true,                                  # Forward
    1,                                      # Wrap (wdFindContinue)
    $false,                                 # Format
    "This is a Microsoft word document.  ", # ReplaceWith
    2                                       # Replace (wdReplaceOne)
) | Out-Null

# 2) Insert a new run with the red annotation right after that text,
# before the paragraph mark, in the same (first) paragraph.
$p = $d.Paragraphs(1)
$r = $p.Range
$r.End = $r.End - 1          # exclude the paragraph mark
$r.Collapse(0)                # collapse to the end (wdCollapseEnd)

$insertStart = $r.Start
$r.InsertAfter("(This is a change – Version for branch alternate)")

# Give just the newly-inserted text its own run with red font color
# (RGB C00000 -> stored as BGR-ordered OLE color 0x0000C0 = 192).
$newRun = $d.Range($insertStart, $r.End)
$newRun.Font.Color = 192
